# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (column F) counts to the 展览 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(9, 6).Value = 122
$ws1.Cells.Item(11, 6).Value = 4591
$ws1.Cells.Item(12, 6).Value = 6809
$ws1.Cells.Item(15, 6).Value = 3546
$ws1.Cells.Item(24, 6).Value = 3693
$ws1.Cells.Item(26, 6).Value = 4077
$ws1.Cells.Item(27, 6).Value = 4078
$ws1.Cells.Item(29, 6).Value = 1932
$ws1.Cells.Item(31, 6).Value = 256
$ws1.Cells.Item(32, 6).Value = 6947
$ws1.Cells.Item(35, 6).Value = 2108
$ws1.Cells.Item(36, 6).Value = 2052
$ws1.Cells.Item(39, 6).Value = 1093
$ws1.Cells.Item(41, 6).Value = 231
$ws1.Cells.Item(43, 6).Value = 227
$ws1.Cells.Item(45, 6).Value = 1151
$ws1.Cells.Item(48, 6).Value = 1850

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 650

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(6, 6).Value = 650
$ws4.Cells.Item(11, 6).Value = 122
$ws4.Cells.Item(14, 6).Value = 4591
$ws4.Cells.Item(15, 6).Value = 6809
$ws4.Cells.Item(17, 6).Value = 3546
$ws4.Cells.Item(28, 6).Value = 4078
$ws4.Cells.Item(31, 6).Value = 256
$ws4.Cells.Item(32, 6).Value = 6947
$ws4.Cells.Item(36, 6).Value = 2108
$ws4.Cells.Item(37, 6).Value = 2052
$ws4.Cells.Item(40, 6).Value = 1093
$ws4.Cells.Item(41, 6).Value = 231
$ws4.Cells.Item(42, 6).Value = 227
$ws4.Cells.Item(44, 6).Value = 1151
$ws4.Cells.Item(47, 6).Value = 1850
